# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets,
# matching the data refresh reflected in the commit (gh-pages data regenerated).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 322
    $ws.Range("F4").Value = 55
    $ws.Range("F5").Value = 280
}
